$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 111; this shifts the previous rows 111..189
# down to 112..190 (and updates the sheet dimension accordingly).
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new record.
$ws.Range("A111").Value2 = 5
$ws.Range("B111").Value2 = "Macroferia Regional de Talca"
$ws.Range("C111").Value2 = "Maule"
$ws.Range("D111").Value2 = 44603
$ws.Range("E111").Value2 = 7
$ws.Range("F111").Value2 = 100112024
$ws.Range("G111").Value2 = "Choclo"
$ws.Range("H111").Value2 = "Choclero"
$ws.Range("I111").Value2 = "Primera"
$ws.Range("J111").Value2 = 40000
$ws.Range("K111").Value2 = 120
$ws.Range("L111").Value2 = 120
$ws.Range("M111").Value2 = 120
$ws.Range("N111").Value2 = "$/unidad"
$ws.Range("O111").Value2 = "Región del Maule"
$ws.Range("P111").Value2 = 120
$ws.Range("Q111").Value2 = 1
$ws.Range("R111").Value2 = "Hortaliza"
